# Generate Report for Handoff
# Adds a new row describing the "ee1cb1bb-4ccc-4b59-a36c-76b692ec4d4c..." file
# (status "Ready for handoff") to the Overview / zh-cn / de-de sheets+tables.

$wb = $excel.ActiveWorkbook

$newMdName   = "ee1cb1bb-4ccc-4b59-a36c-76b692ec4d4cooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newMdE2E    = "e2e\ee1cb1bb-4ccc-4b59-a36c-76b692ec4d4cooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newZhXlf    = "ee1cb1bb-4ccc-4b59-a36c-76b692ec4d4coooooooooooooooooooooooooooooooooooooooo.1a484e86b83b90c43be86fcc9ed625b0632c6a05.zh-cn.xlf"
$newDeXlf    = "ee1cb1bb-4ccc-4b59-a36c-76b692ec4d4coooooooooooooooooooooooooooooooooooooooo.1a484e86b83b90c43be86fcc9ed625b0632c6a05.de-de.xlf"
$newMdUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/277a3b7edacb525c0b2b4c0a326f0200f54bcd85/e2e/$newMdName"

$readyStatus = "Ready for handoff"
$tsHandoff   = "2016-09-02 20:32:56"
$tsZh        = "2016-09-02 20:32:51"
$tsDe        = $tsHandoff
$epoch       = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet (table3) - append summary row
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$rowOv = $loOv.ListRows.Add()
$rOv = $rowOv.Range.Row

$wsOv.Cells.Item($rOv, 1).Value = $newMdName
$wsOv.Cells.Item($rOv, 2).Value = $newMdE2E
$wsOv.Cells.Item($rOv, 3).Value = ".md"
$wsOv.Cells.Item($rOv, 4).Value = ""
$wsOv.Cells.Item($rOv, 5).Value = $readyStatus
$wsOv.Cells.Item($rOv, 6).Value = $readyStatus
$wsOv.Cells.Item($rOv, 7).Value = $tsHandoff

$wsOv.Cells.Item($rOv, 2).Style = "HyperLink"
$wsOv.Hyperlinks.Add($wsOv.Cells.Item($rOv, 2), $newMdUrl, "", "", $newMdE2E)
$wsOv.Cells.Item($rOv, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# zh-cn sheet (table1) - append detail row
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()
$rZh = $rowZh.Range.Row

$wsZh.Cells.Item($rZh, 1).Value = $newMdName
$wsZh.Cells.Item($rZh, 2).Value = ".md"
$wsZh.Cells.Item($rZh, 3).Value = $readyStatus
$wsZh.Cells.Item($rZh, 4).Value = "e2e"
$wsZh.Cells.Item($rZh, 5).Value = "ht"
$wsZh.Cells.Item($rZh, 6).Value = "False"
$wsZh.Cells.Item($rZh, 7).Value = $newZhXlf
$wsZh.Cells.Item($rZh, 8).Value = $tsZh
$wsZh.Cells.Item($rZh, 9).Value = ""
$wsZh.Cells.Item($rZh, 10).Value = ""
$wsZh.Cells.Item($rZh, 11).Value = $epoch
$wsZh.Cells.Item($rZh, 12).Value = ""
$wsZh.Cells.Item($rZh, 13).Value = "True"
$wsZh.Cells.Item($rZh, 14).Value = ""
$wsZh.Cells.Item($rZh, 15).Value = "False"
$wsZh.Cells.Item($rZh, 16).Value = ""

$wsZh.Cells.Item($rZh, 1).Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Cells.Item($rZh, 1), $newMdUrl, "", "", $newMdName)
$wsZh.Cells.Item($rZh, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item($rZh, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# de-de sheet (table2) - append detail row
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()
$rDe = $rowDe.Range.Row

$wsDe.Cells.Item($rDe, 1).Value = $newMdName
$wsDe.Cells.Item($rDe, 2).Value = ".md"
$wsDe.Cells.Item($rDe, 3).Value = $readyStatus
$wsDe.Cells.Item($rDe, 4).Value = "e2e"
$wsDe.Cells.Item($rDe, 5).Value = "ht"
$wsDe.Cells.Item($rDe, 6).Value = "False"
$wsDe.Cells.Item($rDe, 7).Value = $newDeXlf
$wsDe.Cells.Item($rDe, 8).Value = $tsDe
$wsDe.Cells.Item($rDe, 9).Value = ""
$wsDe.Cells.Item($rDe, 10).Value = ""
$wsDe.Cells.Item($rDe, 11).Value = $epoch
$wsDe.Cells.Item($rDe, 12).Value = ""
$wsDe.Cells.Item($rDe, 13).Value = "True"
$wsDe.Cells.Item($rDe, 14).Value = ""
$wsDe.Cells.Item($rDe, 15).Value = "False"
$wsDe.Cells.Item($rDe, 16).Value = ""

$wsDe.Cells.Item($rDe, 1).Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Cells.Item($rDe, 1), $newMdUrl, "", "", $newMdName)
$wsDe.Cells.Item($rDe, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item($rDe, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Column width tweaks (observed in diff: widened date columns)
# ---------------------------------------------------------------------------
$wsOv.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOv.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZh.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDe.Columns.Item(3).ColumnWidth = 17.2159881591797

Write-Host "done"
